$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 target values
$ws.Range("A2").Value = 111809580
$ws.Range("B2").Value = 77515
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 610571.4165256479
$ws.Range("R2").Value = 7180702.680798599
$ws.Range("Z2").Value = "14:31"
$ws.Range("AB2").Value = "14:31"

# Row 3 target values
$ws.Range("A3").Value = 111809638
$ws.Range("B3").Value = 56398
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("Q3").Value = 610564.5541715119
$ws.Range("R3").Value = 7180691.309759256
$ws.Range("Z3").Value = "14:36"
$ws.Range("AB3").Value = "14:36"

# Row 4 target values
$ws.Range("A4").Value = 111809656
$ws.Range("B4").Value = 56398
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 610542.5625081829
$ws.Range("R4").Value = 7180707.182562917
$ws.Range("Z4").Value = "14:37"
$ws.Range("AB4").Value = "14:37"
